# Updated cryptos list on Thu May  4 22:50:29 UTC 2023 with GitHub Actions
#
# Refresh the Price (D) and Volume(1h) (E) columns for the crypto table on
# the active sheet. Values are written as plain text (matching the source
# data, which stores prices/percentages as strings, not numbers) - a
# leading '' (an escaped single quote inside a single-quoted PowerShell
# string) forces Excel to keep numeric-looking price strings as literal
# text instead of auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.857.87'
$ws.Range("E2").Value = '  -1.42%  '
$ws.Range("D3").Value = '1.877.16'
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '''324.53'
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("D6").Value = '''1.003'
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").Value = '''0.4612'
$ws.Range("E7").Value = '  -1.30%  '
$ws.Range("D8").Value = '''0.3872'
$ws.Range("E8").Value = '  -2.59%  '
$ws.Range("D9").Value = '''0.07844'
$ws.Range("E9").Value = '  -2.59%  '
$ws.Range("D10").Value = '''0.9824'
$ws.Range("D11").Value = '''21.71'
$ws.Range("E11").Value = '  -2.70%  '
$ws.Range("D12").Value = '1.895.19'
$ws.Range("E12").Value = '  -0.47%  '
$ws.Range("D13").Value = '''6.988'
$ws.Range("E13").Value = '  -2.39%  '
$ws.Range("D14").Value = '''5.662'
$ws.Range("E14").Value = '  -2.71%  '
$ws.Range("D15").Value = '''0.06966'
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D16").Value = '''88.18'
$ws.Range("E16").Value = '  -1.78%  '
$ws.Range("D18").Value = '''0.000009950'
$ws.Range("E18").Value = '  -2.90%  '
$ws.Range("D19").Value = '''16.89'
$ws.Range("E19").Value = '  -2.92%  '
$ws.Range("D20").Value = '''1.002'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").Value = '28.857.61'
$ws.Range("E21").Value = '  -1.47%  '
$ws.Range("D22").Value = '''5.259'
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").Value = '''10.98'
$ws.Range("E23").Value = '  -1.87%  '
$ws.Range("D24").Value = '''2.106'
$ws.Range("E24").Value = '  +1.96%  '
$ws.Range("D25").Value = '''156.12'
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("D26").Value = '''19.36'
$ws.Range("E26").Value = '  -2.12%  '
$ws.Range("D27").Value = '''5.946'
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").Value = '''117.63'
$ws.Range("E28").Value = '  -2.87%  '
$ws.Range("D29").Value = '''1.904'
$ws.Range("E29").Value = '  -6.52%  '
$ws.Range("D30").Value = '''0.09352'
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("D31").Value = '''0.8991'
$ws.Range("E31").Value = '  -4.92%  '
$ws.Range("D32").Value = '''5.263'
$ws.Range("E32").Value = '  -2.44%  '
$ws.Range("E33").Value = '  -2.84%  '
$ws.Range("D34").Value = '''3.250'
$ws.Range("E34").Value = '  -0.64%  '
$ws.Range("D35").Value = '''1.170'
$ws.Range("E35").Value = '  -0.92%  '
$ws.Range("D36").Value = '''0.05732'
$ws.Range("E36").Value = '  -2.80%  '
$ws.Range("D37").Value = '''0.02071'
$ws.Range("E37").Value = '  -2.13%  '
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("D39").Value = '''7.617'
$ws.Range("E39").Value = '  -6.28%  '
$ws.Range("D40").Value = '''0.5651'
$ws.Range("E40").Value = '  -3.48%  '
$ws.Range("D41").Value = '''0.1770'
$ws.Range("E41").Value = '  -2.93%  '
$ws.Range("D42").Value = '''9.690'
$ws.Range("E42").Value = '  -4.36%  '
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("D44").Value = '''2.232'
$ws.Range("E44").Value = '  -4.06%  '
$ws.Range("D45").Value = '''0.5330'
$ws.Range("E45").Value = '  -2.84%  '
$ws.Range("D46").Value = '''0.07040'
$ws.Range("E46").Value = '  -2.69%  '
$ws.Range("D47").Value = '''1.839'
$ws.Range("E47").Value = '  -4.25%  '
$ws.Range("D48").Value = '''2.538'
$ws.Range("E48").Value = '  +0.60%  '
$ws.Range("D49").Value = '''112.37'
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("D50").Value = '''1.053'
$ws.Range("E50").Value = '  -7.16%  '
$ws.Range("D51").Value = '''70.77'
$ws.Range("E51").Value = '  -1.54%  '
